# Add gid and lang columns (data source now reflects a DB connection providing
# a global id and language for each localized scenario/title row).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "title" column (B), which will
# become columns B (gid) and C (lang); the old "title" column shifts to D.
$ws.Range("B1:C1").EntireColumn.Insert()

# Headers (D1 already holds "title", shifted over by the column insert above)
$ws.Range("B1").Value = "gid"
$ws.Range("C1").Value = "lang"

# Keep gid values (large numeric ids) stored as text, matching the source data.
$ws.Range("B2:B14").NumberFormat = "@"

$ws.Range("B2").Value = "3095466740"
$ws.Range("C2").Value = "EN"
$ws.Range("B3").Value = "1019682570"
$ws.Range("C3").Value = "EN"
$ws.Range("B4").Value = "2879327918"
$ws.Range("C4").Value = "EN"
$ws.Range("B5").Value = "1042118304"
$ws.Range("C5").Value = "EN"
$ws.Range("B6").Value = "1323093475"
$ws.Range("C6").Value = "EN"
$ws.Range("B7").Value = "3023688674"
$ws.Range("C7").Value = "EN"
$ws.Range("B8").Value = "2322922310"
$ws.Range("C8").Value = "EN"
$ws.Range("B9").Value = "1400916808"
$ws.Range("C9").Value = "EN"
$ws.Range("B10").Value = "2204235069"
$ws.Range("C10").Value = "EN"
$ws.Range("B11").Value = "2038411104"
$ws.Range("C11").Value = "EN"
$ws.Range("B12").Value = "250690944"
$ws.Range("C12").Value = "EN"
$ws.Range("B13").Value = "1897449103"
$ws.Range("C13").Value = "EN"
$ws.Range("B14").Value = "1997989660"
$ws.Range("C14").Value = "EN"
